# Insert a new "statut_name" column before the existing NCTId column (C),
# shifting NCTId..intervention_type from C:L to D:M, then populate the new
# column with a human-readable status label derived from the
# results_1y / results_3y / results booleans (now shifted to J/K/L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at C; existing C:L -> D:M.
$ws.Columns.Item(3).Insert()

# Header for the new column.
$ws.Cells.Item(1, 3).Value = "statut_name"

$lastRow = 74

for ($r = 2; $r -le $lastRow; $r++) {
    $results1y = $ws.Cells.Item($r, 10).Text   # was column I (results_1y), now J
    $results3y = $ws.Cells.Item($r, 11).Text   # was column J (results_3y), now K
    $results   = $ws.Cells.Item($r, 12).Text   # was column K (results),    now L

    if ($results1y -eq "TRUE") {
        $name = "résultat et / ou publication posté dans les 12 mois"
    } elseif ($results3y -eq "TRUE") {
        $name = "résultat et / ou publication posté dans les 36 mois"
    } elseif ($results -eq "TRUE") {
        $name = "résultat et / ou publication posté"
    } else {
        $name = "pas de résultat ni de publication"
    }

    $ws.Cells.Item($r, 3).Value = $name
}
